# Update to Jogos_da_Semana_FlashScore_2024-10-24.xlsx
# 1) Tweak several odds in row 5 (Sport Recife - Guarani)
# 2) Insert a brand-new match row at position 11 (Correcaminos - Tapatio,
#    Mexico - Liga de Expansion MX), pushing the former rows 11-19 down to 12-20

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Row 5 odds updates -------------------------------------------------
$row5Updates = @{
    "G5"  = 1.57
    "I5"  = 7
    "K5"  = 1.95
    "L5"  = 8
    "U5"  = 3
    "V5"  = 1.36
    "X5"  = 5.5
    "AD5" = 8
    "AG5" = 11
    "AH5" = 34
    "AN5" = 3.2
    "AO5" = 8.5
    "AQ5" = 29
    "AZ5" = 251
    "BA5" = 351
}

foreach ($addr in $row5Updates.Keys) {
    $ws.Range($addr).Value = $row5Updates[$addr]
}

# --- 2) Insert new row 11 and shift existing rows 11-19 down to 12-20 ------
$ws.Range("A11").EntireRow.Insert()

# Values for the newly inserted row 11, in column order A..BD
$newRow11 = @(
    "8QmVgVD5",
    "24/10/2024",
    "22:00",
    "MEXICO - LIGA DE EXPANSION MX",
    "Correcaminos",
    "Tapatio",
    3.75,
    3.5,
    1.85,
    4.15,
    2.18,
    2.4,
    1.01,
    11,
    1.23,
    3.35,
    1.7,
    1.93,
    1.34,
    3.1,
    1.62,
    2.02,
    12,
    22,
    12.5,
    55,
    32,
    35,
    11.5,
    6.9,
    13.5,
    55,
    8.25,
    9.75,
    8.25,
    16.5,
    14,
    22,
    350,
    5.7,
    21,
    25,
    110,
    120,
    300,
    2.9,
    6.8,
    55,
    3.8,
    9,
    16,
    30,
    55,
    175,
    51,
    51
)

for ($i = 0; $i -lt $newRow11.Length; $i++) {
    $ws.Cells.Item(11, $i + 1).Value = $newRow11[$i]
}
